$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the base frequency value (B1): 400 -> 6400
$ws.Range("B1").Value = 6400

# Move the selection from B1 to B2
$ws.Range("B2").Select()
